# Update the WRESBAL FRED data sheet with newly published observations.
#
# The source data rolls forward each time it is refreshed:
#   - 5 older weekly observations (Sep 2021) are now available and get
#     inserted at the top of the table (right after the header row).
#   - 2 newer weekly observations (Nov/Dec 2023) are appended at the end.
#   - The SeriesInfo sheet metadata (realtime_start/end, observation_end,
#     last_updated) is refreshed to reflect the new pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1) Insert 5 new rows right after the header (before the old first
#    data row) and format them like the existing date/value rows.
# ---------------------------------------------------------------------
$ws.Rows("2:6").Insert()

# Copy the formatting from the (now shifted) first original data row
# - row 7, which used to be row 2 - onto the newly inserted blank rows.
$ws.Range("A7:B7").Copy()
$ws.Range("A2:B6").PasteSpecial(-4122)

$ws.Range("A2").Value = 44440
$ws.Range("B2").Value = 4193.981

$ws.Range("A3").Value = 44447
$ws.Range("B3").Value = 4250.506

$ws.Range("A4").Value = 44454
$ws.Range("B4").Value = 4271.184

$ws.Range("A5").Value = 44461
$ws.Range("B5").Value = 4156.747

$ws.Range("A6").Value = 44468
$ws.Range("B6").Value = 4125.348

# ---------------------------------------------------------------------
# 2) Append 2 new rows at the end of the table (after the original
#    last row, now at row 118) with the newest observations.
# ---------------------------------------------------------------------
$ws.Range("A118:B118").Copy()
$ws.Range("A119:B120").PasteSpecial(-4122)

$ws.Range("A119").Value = 45259
$ws.Range("B119").Value = 3441.581

$ws.Range("A120").Value = 45266
$ws.Range("B120").Value = 3513.374

# ---------------------------------------------------------------------
# 3) Refresh the SeriesInfo metadata sheet.
#    NOTE: B3/B4/B7 hold plain "YYYY-MM-DD" text in the source file, but
#    Excel's smart-entry would otherwise reinterpret that pattern as a
#    date serial. A leading apostrophe forces literal text, exactly like
#    typing it into the Excel UI would.
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

$info.Range("B3").Value = "'2023-12-08"
$info.Range("B4").Value = "'2023-12-08"
$info.Range("B7").Value = "'2023-12-06"
$info.Range("B14").Value = "2023-12-07 15:35:02-06"
